$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47 (shifts old rows 47-99 down to 48-100)
$ws.Rows.Item(47).Insert()

# Write all rows from 47 to 142 with final content
$ws.Range("A47").Value = 'blog.metadata.separator'
$ws.Range("B47").Value = '•'
$ws.Range("C47").Value = '•'

$ws.Range("A48").Value = 'resources.title'
$ws.Range("B48").Value = 'Resources'
$ws.Range("C48").Value = 'Tài Nguyên'

$ws.Range("A49").Value = 'resources.description'
$ws.Range("B49").Value = 'A collection of tools, guides, and resources to help you on your journey'
$ws.Range("C49").Value = 'Bộ sưu tập công cụ, hướng dẫn và tài nguyên để hỗ trợ bạn trên hành trình của mình'

$ws.Range("A50").Value = 'resources.tools'
$ws.Range("B50").Value = 'Tools'
$ws.Range("C50").Value = 'Công Cụ'

$ws.Range("A51").Value = 'resources.guides'
$ws.Range("B51").Value = 'Guides'
$ws.Range("C51").Value = 'Hướng Dẫn'

$ws.Range("A52").Value = 'about.title'
$ws.Range("B52").Value = 'About Me'
$ws.Range("C52").Value = 'Về Tôi'

$ws.Range("A53").Value = 'about.intro'
$ws.Range("B53").Value = 'Hi, I''m Jonathan. I''m passionate about technology and making a positive impact.'
$ws.Range("C53").Value = 'Xin chào, tôi là Jonathan. Tôi đam mê công nghệ và tạo ra những tác động tích cực.'

$ws.Range("A54").Value = 'about.mission'
$ws.Range("B54").Value = 'My mission is to help others leverage technology for positive change.'
$ws.Range("C54").Value = 'Sứ mệnh của tôi là giúp mọi người tận dụng công nghệ để tạo ra những thay đổi tích cực.'

$ws.Range("A55").Value = 'about.contact'
$ws.Range("B55").Value = 'Get in Touch'
$ws.Range("C55").Value = 'Liên Hệ'

$ws.Range("A56").Value = 'about.background.title'
$ws.Range("B56").Value = 'Background'
$ws.Range("C56").Value = 'Giới Thiệu'

$ws.Range("A57").Value = 'about.background.content'
$ws.Range("B57").Value = 'As a Product Manager with a passion for technology and social impact, I bridge the gap between innovative solutions and human needs. My journey in product management has been driven by a commitment to creating meaningful digital experiences that make a difference.'
$ws.Range("C57").Value = 'Là một Product Manager với niềm đam mê về công nghệ và tác động xã hội, tôi kết nối giữa các giải pháp sáng tạo và nhu cầu con người. Hành trình của tôi trong quản lý sản phẩm được thúc đẩy bởi cam kết tạo ra những trải nghiệm kỹ thuật số có ý nghĩa.'

$ws.Range("A58").Value = 'about.expertise.title'
$ws.Range("B58").Value = 'Expertise'
$ws.Range("C58").Value = 'Chuyên Môn'

$ws.Range("A59").Value = 'about.expertise.skills.product_strategy'
$ws.Range("B59").Value = 'Product Strategy'
$ws.Range("C59").Value = 'Chiến Lược Sản Phẩm'

$ws.Range("A60").Value = 'about.expertise.skills.user_research'
$ws.Range("B60").Value = 'User Research'
$ws.Range("C60").Value = 'Nghiên Cứu Người Dùng'

$ws.Range("A61").Value = 'about.expertise.skills.agile'
$ws.Range("B61").Value = 'Agile Management'
$ws.Range("C61").Value = 'Quản Lý Agile'

$ws.Range("A62").Value = 'about.expertise.skills.analytics'
$ws.Range("B62").Value = 'Data Analytics'
$ws.Range("C62").Value = 'Phân Tích Dữ Liệu'

$ws.Range("A63").Value = 'about.expertise.skills.technical'
$ws.Range("B63").Value = 'Technical Leadership'
$ws.Range("C63").Value = 'Lãnh Đạo Kỹ Thuật'

$ws.Range("A64").Value = 'about.vision.title'
$ws.Range("B64").Value = 'Vision'
$ws.Range("C64").Value = 'Tầm Nhìn'

$ws.Range("A65").Value = 'about.vision.content'
$ws.Range("B65").Value = 'I believe in leveraging technology to create products that not only solve problems but also contribute positively to society. My goal is to lead product initiatives that combine innovation with social responsibility.'
$ws.Range("C65").Value = 'Tôi tin vào việc tận dụng công nghệ để tạo ra những sản phẩm không chỉ giải quyết vấn đề mà còn đóng góp tích cực cho xã hội. Mục tiêu của tôi là dẫn dắt các sáng kiến sản phẩm kết hợp đổi mới với trách nhiệm xã hội.'

$ws.Range("A66").Value = 'about.timeline.present.year'
$ws.Range("B66").Value = '2023 - Present'
$ws.Range("C66").Value = '2023 - Hiện tại'

$ws.Range("A67").Value = 'about.timeline.present.title'
$ws.Range("B67").Value = 'Senior Product Manager'
$ws.Range("C67").Value = 'Quản Lý Sản Phẩm Cao Cấp'

$ws.Range("A68").Value = 'about.timeline.present.description'
$ws.Range("B68").Value = 'Leading innovative product initiatives'
$ws.Range("C68").Value = 'Dẫn dắt các sáng kiến sản phẩm sáng tạo'

$ws.Range("A69").Value = 'about.timeline.past_1.year'
$ws.Range("B69").Value = '2020 - 2023'
$ws.Range("C69").Value = '2020 - 2023'

$ws.Range("A70").Value = 'about.timeline.past_1.title'
$ws.Range("B70").Value = 'Product Manager'
$ws.Range("C70").Value = 'Quản Lý Sản Phẩm'

$ws.Range("A71").Value = 'about.timeline.past_1.description'
$ws.Range("B71").Value = 'Driving user-centered product development'
$ws.Range("C71").Value = 'Thúc đẩy phát triển sản phẩm lấy người dùng làm trung tâm'

$ws.Range("A72").Value = 'about.timeline.past_2.year'
$ws.Range("B72").Value = '2018 - 2020'
$ws.Range("C72").Value = '2018 - 2020'

$ws.Range("A73").Value = 'about.timeline.past_2.title'
$ws.Range("B73").Value = 'Associate Product Manager'
$ws.Range("C73").Value = 'Quản Lý Sản Phẩm Phó'

$ws.Range("A74").Value = 'about.timeline.past_2.description'
$ws.Range("B74").Value = 'Building foundation in product management'
$ws.Range("C74").Value = 'Xây dựng nền tảng trong quản lý sản phẩm'

$ws.Range("A75").Value = 'contact.title'
$ws.Range("B75").Value = 'Get in Touch'
$ws.Range("C75").Value = 'Get in Touch'

$ws.Range("A76").Value = 'contact.description'
$ws.Range("B76").Value = 'Have a question or want to collaborate? I''d love to hear from you.'
$ws.Range("C76").Value = 'Have a question or want to collaborate? I''d love to hear from you.'

$ws.Range("A77").Value = 'contact.form.name'
$ws.Range("B77").Value = 'Your Name'
$ws.Range("C77").Value = 'Your Name'

$ws.Range("A78").Value = 'contact.form.email'
$ws.Range("B78").Value = 'Your Email'
$ws.Range("C78").Value = 'Your Email'

$ws.Range("A79").Value = 'contact.form.message'
$ws.Range("B79").Value = 'Your Message'
$ws.Range("C79").Value = 'Your Message'

$ws.Range("A80").Value = 'contact.form.submit'
$ws.Range("B80").Value = 'Send Message'
$ws.Range("C80").Value = 'Send Message'

$ws.Range("A81").Value = 'contact.form.success'
$ws.Range("B81").Value = 'Message sent successfully!'
$ws.Range("C81").Value = 'Message sent successfully!'

$ws.Range("A82").Value = 'contact.form.error'
$ws.Range("B82").Value = 'Error sending message. Please try again.'
$ws.Range("C82").Value = 'Error sending message. Please try again.'

$ws.Range("A83").Value = 'finance.title'
$ws.Range("B83").Value = 'Personal Finance'
$ws.Range("C83").Value = 'Quản Lý Tài Chính Cá Nhân'

$ws.Range("A84").Value = 'finance.description'
$ws.Range("B84").Value = 'Insights and strategies for building financial independence'
$ws.Range("C84").Value = 'Những hiểu biết và chiến lược để xây dựng tự do tài chính'

$ws.Range("A85").Value = 'finance.sections.investments.title'
$ws.Range("B85").Value = 'Investment Strategies'
$ws.Range("C85").Value = 'Chiến Lược Đầu Tư'

$ws.Range("A86").Value = 'finance.sections.investments.description'
$ws.Range("B86").Value = 'Long-term approaches to wealth building'
$ws.Range("C86").Value = 'Phương pháp xây dựng tài sản dài hạn'

$ws.Range("A87").Value = 'finance.sections.budgeting.title'
$ws.Range("B87").Value = 'Smart Budgeting'
$ws.Range("C87").Value = 'Quản Lý Chi Tiêu'

$ws.Range("A88").Value = 'finance.sections.budgeting.description'
$ws.Range("B88").Value = 'Practical tips for effective money management'
$ws.Range("C88").Value = 'Mẹo thực tế để quản lý tiền hiệu quả'

$ws.Range("A89").Value = 'finance.sections.planning.title'
$ws.Range("B89").Value = 'Financial Planning'
$ws.Range("C89").Value = 'Kế Hoạch Tài Chính'

$ws.Range("A90").Value = 'finance.sections.planning.description'
$ws.Range("B90").Value = 'Setting and achieving financial goals'
$ws.Range("C90").Value = 'Thiết lập và đạt được mục tiêu tài chính'

$ws.Range("A91").Value = 'impact.title'
$ws.Range("B91").Value = 'Social Impact'
$ws.Range("C91").Value = 'Tác Động Xã Hội'

$ws.Range("A92").Value = 'impact.description'
$ws.Range("B92").Value = 'Making a positive difference in the world'
$ws.Range("C92").Value = 'Tạo ra những thay đổi tích cực cho thế giới'

$ws.Range("A93").Value = 'impact.sections.projects.title'
$ws.Range("B93").Value = 'Impact Projects'
$ws.Range("C93").Value = 'Dự Án Tác Động'

$ws.Range("A94").Value = 'impact.sections.projects.description'
$ws.Range("B94").Value = 'Current initiatives and their outcomes'
$ws.Range("C94").Value = 'Các sáng kiến hiện tại và kết quả'

$ws.Range("A95").Value = 'impact.sections.metrics.title'
$ws.Range("B95").Value = 'Impact Metrics'
$ws.Range("C95").Value = 'Đo Lường Tác Động'

$ws.Range("A96").Value = 'impact.sections.metrics.description'
$ws.Range("B96").Value = 'Measuring and tracking social impact'
$ws.Range("C96").Value = 'Đo lường và theo dõi tác động xã hội'

$ws.Range("A97").Value = 'impact.sections.collaboration.title'
$ws.Range("B97").Value = 'Get Involved'
$ws.Range("C97").Value = 'Tham Gia'

$ws.Range("A98").Value = 'impact.sections.collaboration.description'
$ws.Range("B98").Value = 'Ways to contribute and collaborate'
$ws.Range("C98").Value = 'Cách đóng góp và hợp tác'

$ws.Range("A99").Value = 'common.back_to_library'
$ws.Range("B99").Value = 'Back to Library'
$ws.Range("C99").Value = 'Quay Lại Thư Viện'

$ws.Range("A100").Value = 'common.published_on'
$ws.Range("B100").Value = 'Published on {{date}}'
$ws.Range("C100").Value = 'Xuất bản ngày {{date}}'

$ws.Range("A101").Value = 'common.explore'
$ws.Range("B101").Value = 'Explore Section'
$ws.Range("C101").Value = 'Khám Phá'

$ws.Range("A102").Value = 'auth.login.title'
$ws.Range("B102").Value = 'Sign in to your account'
$ws.Range("C102").Value = 'Đăng nhập vào tài khoản của bạn'

$ws.Range("A103").Value = 'auth.login.email_label'
$ws.Range("B103").Value = 'Email address'
$ws.Range("C103").Value = 'Địa chỉ email'

$ws.Range("A104").Value = 'auth.login.email_placeholder'
$ws.Range("B104").Value = 'Email address'
$ws.Range("C104").Value = 'Địa chỉ email'

$ws.Range("A105").Value = 'auth.login.password_label'
$ws.Range("B105").Value = 'Password'
$ws.Range("C105").Value = 'Mật khẩu'

$ws.Range("A106").Value = 'auth.login.password_placeholder'
$ws.Range("B106").Value = 'Password'
$ws.Range("C106").Value = 'Mật khẩu'

$ws.Range("A107").Value = 'auth.login.forgot_password'
$ws.Range("B107").Value = 'Forgot your password?'
$ws.Range("C107").Value = 'Quên mật khẩu?'

$ws.Range("A108").Value = 'auth.login.submit_button'
$ws.Range("B108").Value = 'Sign in'
$ws.Range("C108").Value = 'Đăng nhập'

$ws.Range("A109").Value = 'auth.login.submitting'
$ws.Range("B109").Value = 'Signing in...'
$ws.Range("C109").Value = 'Đang đăng nhập...'

$ws.Range("A110").Value = 'auth.login.need_account'
$ws.Range("B110").Value = 'Need an account?'
$ws.Range("C110").Value = 'Chưa có tài khoản?'

$ws.Range("A111").Value = 'auth.login.sign_up_link'
$ws.Range("B111").Value = 'Sign up'
$ws.Range("C111").Value = 'Đăng ký'

$ws.Range("A112").Value = 'auth.login.error'
$ws.Range("B112").Value = 'Failed to log in: {{message}}'
$ws.Range("C112").Value = 'Đăng nhập thất bại: {{message}}'

$ws.Range("A113").Value = 'auth.signup.title'
$ws.Range("B113").Value = 'Create your account'
$ws.Range("C113").Value = 'Tạo tài khoản của bạn'

$ws.Range("A114").Value = 'auth.signup.email_label'
$ws.Range("B114").Value = 'Email address'
$ws.Range("C114").Value = 'Địa chỉ email'

$ws.Range("A115").Value = 'auth.signup.email_placeholder'
$ws.Range("B115").Value = 'Email address'
$ws.Range("C115").Value = 'Địa chỉ email'

$ws.Range("A116").Value = 'auth.signup.password_label'
$ws.Range("B116").Value = 'Password'
$ws.Range("C116").Value = 'Mật khẩu'

$ws.Range("A117").Value = 'auth.signup.password_placeholder'
$ws.Range("B117").Value = 'Password'
$ws.Range("C117").Value = 'Mật khẩu'

$ws.Range("A118").Value = 'auth.signup.confirm_password_label'
$ws.Range("B118").Value = 'Confirm Password'
$ws.Range("C118").Value = 'Xác nhận mật khẩu'

$ws.Range("A119").Value = 'auth.signup.confirm_password_placeholder'
$ws.Range("B119").Value = 'Confirm Password'
$ws.Range("C119").Value = 'Xác nhận mật khẩu'

$ws.Range("A120").Value = 'auth.signup.submit_button'
$ws.Range("B120").Value = 'Sign up'
$ws.Range("C120").Value = 'Đăng ký'

$ws.Range("A121").Value = 'auth.signup.submitting'
$ws.Range("B121").Value = 'Creating account...'
$ws.Range("C121").Value = 'Đang tạo tài khoản...'

$ws.Range("A122").Value = 'auth.signup.have_account'
$ws.Range("B122").Value = 'Already have an account?'
$ws.Range("C122").Value = 'Đã có tài khoản?'

$ws.Range("A123").Value = 'auth.signup.login_link'
$ws.Range("B123").Value = 'Log in'
$ws.Range("C123").Value = 'Đăng nhập'

$ws.Range("A124").Value = 'auth.signup.error'
$ws.Range("B124").Value = 'Failed to create an account: {{message}}'
$ws.Range("C124").Value = 'Tạo tài khoản thất bại: {{message}}'

$ws.Range("A125").Value = 'auth.signup.passwords_not_match'
$ws.Range("B125").Value = 'Passwords do not match'
$ws.Range("C125").Value = 'Mật khẩu không khớp'

$ws.Range("A126").Value = 'auth.forgot_password.title'
$ws.Range("B126").Value = 'Reset your password'
$ws.Range("C126").Value = 'Đặt lại mật khẩu của bạn'

$ws.Range("A127").Value = 'auth.forgot_password.email_label'
$ws.Range("B127").Value = 'Email address'
$ws.Range("C127").Value = 'Địa chỉ email'

$ws.Range("A128").Value = 'auth.forgot_password.email_placeholder'
$ws.Range("B128").Value = 'Email address'
$ws.Range("C128").Value = 'Địa chỉ email'

$ws.Range("A129").Value = 'auth.forgot_password.submit_button'
$ws.Range("B129").Value = 'Reset Password'
$ws.Range("C129").Value = 'Đặt lại mật khẩu'

$ws.Range("A130").Value = 'auth.forgot_password.submitting'
$ws.Range("B130").Value = 'Processing...'
$ws.Range("C130").Value = 'Đang xử lý...'

$ws.Range("A131").Value = 'auth.forgot_password.back_to_login'
$ws.Range("B131").Value = 'Back to login'
$ws.Range("C131").Value = 'Quay lại đăng nhập'

$ws.Range("A132").Value = 'auth.forgot_password.error'
$ws.Range("B132").Value = 'Failed to reset password: {{message}}'
$ws.Range("C132").Value = 'Đặt lại mật khẩu thất bại: {{message}}'

$ws.Range("A133").Value = 'auth.forgot_password.success_message'
$ws.Range("B133").Value = 'Check your inbox for further instructions'
$ws.Range("C133").Value = 'Kiểm tra hộp thư của bạn để biết hướng dẫn tiếp theo'

$ws.Range("A134").Value = 'auth.profile.title'
$ws.Range("B134").Value = 'Profile'
$ws.Range("C134").Value = 'Hồ sơ'

$ws.Range("A135").Value = 'auth.profile.user_info_title'
$ws.Range("B135").Value = 'User Information'
$ws.Range("C135").Value = 'Thông tin người dùng'

$ws.Range("A136").Value = 'auth.profile.user_info_description'
$ws.Range("B136").Value = 'Personal details and account settings.'
$ws.Range("C136").Value = 'Thông tin cá nhân và cài đặt tài khoản.'

$ws.Range("A137").Value = 'auth.profile.email_label'
$ws.Range("B137").Value = 'Email'
$ws.Range("C137").Value = 'Email'

$ws.Range("A138").Value = 'auth.profile.email_verified_label'
$ws.Range("B138").Value = 'Email verified'
$ws.Range("C138").Value = 'Email đã xác minh'

$ws.Range("A139").Value = 'auth.profile.yes'
$ws.Range("B139").Value = 'Yes'
$ws.Range("C139").Value = 'Có'

$ws.Range("A140").Value = 'auth.profile.no'
$ws.Range("B140").Value = 'No'
$ws.Range("C140").Value = 'Không'

$ws.Range("A141").Value = 'auth.profile.logout_button'
$ws.Range("B141").Value = 'Log Out'
$ws.Range("C141").Value = 'Đăng xuất'

$ws.Range("A142").Value = 'auth.profile.error'
$ws.Range("B142").Value = 'Failed to log out'
$ws.Range("C142").Value = 'Đăng xuất thất bại'
